$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data (including two coin-row swaps) from latest scrape
# Force text number format so numeric-looking strings (e.g. "1.00", "503.77") are
# preserved exactly as text, then restore the default "Normal" style so no stray
# cell formatting is introduced.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '54.690.09'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  +0.61%  '
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.278.97'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.13%  '
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '503.77'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  +1.08%  '
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '128.46'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  +0.11%  '
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.12%  '
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.528'
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -0.02%  '
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '2.294.47'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  +0.36%  '
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0967'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  +0.79%  '
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +0.55%  '
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.343'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  +3.36%  '
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.90'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  +3.99%  '
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '23.39'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  +4.50%  '
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.685.42'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '54.778.79'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  +0.91%  '
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  +0.92%  '
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.287.40'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  +0.48%  '
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '10.33'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  +1.53%  '
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '4.13'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '306.84'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  +0.50%  '
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.42'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -0.59%  '
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -0.30%  '
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '60.15'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -2.88%  '
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  +0.05%  '
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  +0.10%  '
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '7.42'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  +1.67%  '
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '170.83'
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -0.73%  '
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.0₃0701'
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  +1.89%  '
$c.Style = "Normal"

$c = $ws.Range("B30")
$c.NumberFormat = "@"
$c.Value = 'PancakeSwap'
$c.Style = "Normal"
$c = $ws.Range("C30")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.63'
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  +0.48%  '
$c.Style = "Normal"

$c = $ws.Range("B31")
$c.NumberFormat = "@"
$c.Value = 'Aptos'
$c.Style = "Normal"
$c = $ws.Range("C31")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.05'
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  +1.90%  '
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  +2.80%  '
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '17.92'
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +0.87%  '
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.995'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -0.11%  '
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.913'
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -1.27%  '
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.20'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -0.39%  '
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.79'
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  +1.57%  '
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '36.45'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  +1.87%  '
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.375'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  +0.56%  '
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.41'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +0.05%  '
$c.Style = "Normal"

$c = $ws.Range("B42")
$c.NumberFormat = "@"
$c.Value = 'RenderToken'
$c.Style = "Normal"
$c = $ws.Range("C42")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.03'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  +0.75%  '
$c.Style = "Normal"

$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = 'Filecoin'
$c.Style = "Normal"
$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.39'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '126.40'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '252.34'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  +5.42%  '
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0497'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  +1.29%  '
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.0900'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +0.17%  '
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.548'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.374'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  +0.66%  '
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0206'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '10.82'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  +0.46%  '
$c.Style = "Normal"
